# Update Name of Algo
# Apply updated numeric values produced by the (renamed) KNN imputation algorithm.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D3").Value  = -7.812
$ws.Range("D4").Value  = -8.032
$ws.Range("E6").Value  = 12.718
$ws.Range("D7").Value  = -8.055999999999999
$ws.Range("E7").Value  = 12.708
$ws.Range("D8").Value  = -8.010999999999999
$ws.Range("E8").Value  = 13.03
$ws.Range("C11").Value = -12.634
$ws.Range("C12").Value = -12.976
$ws.Range("D12").Value = -7.992999999999999
$ws.Range("D14").Value = -8.224
$ws.Range("C15").Value = -12.493
$ws.Range("E19").Value = 12.36
$ws.Range("E21").Value = 13.136
$ws.Range("D22").Value = -7.811999999999999
$ws.Range("E24").Value = 12.623
$ws.Range("E25").Value = 12.36
